$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 9.462749333333333
$ws.Cells.Item(2, 8).Value = 28.388248
$ws.Cells.Item(2, 9).Value = 0.07254428564686972
$ws.Cells.Item(2, 10).Value = 0.07439525120506714
$ws.Cells.Item(2, 13).Value = 3.790335666666667
$ws.Cells.Item(2, 14).Value = 11.371007
$ws.Cells.Item(2, 15).Value = 0.02257417106507704
$ws.Cells.Item(2, 16).Value = 0.02284136466268864
$ws.Cells.Item(2, 17).Value = 35.86699630285955
$ws.Cells.Item(2, 18).Value = 322.802966725736
$ws.Cells.Item(2, 19).Value = 0.001637627113986251
$ws.Cells.Item(2, 20).Value = 0.001699289061947265
$ws.Cells.Item(3, 7).Value = 9.462749333333333
$ws.Cells.Item(3, 8).Value = 28.388248
$ws.Cells.Item(3, 9).Value = 0.07254428564686972
$ws.Cells.Item(3, 10).Value = 0.07439525120506714
$ws.Cells.Item(3, 15).Value = 0.2106472625291773
$ws.Cells.Item(3, 16).Value = 0.2131405368000222
$ws.Cells.Item(3, 17).Value = 334.6871326774764
$ws.Cells.Item(3, 18).Value = 3012.184194097288
$ws.Cells.Item(3, 19).Value = 0.0152812551836478
$ws.Cells.Item(3, 20).Value = 0.0158566437772205
$ws.Cells.Item(4, 7).Value = 9.462749333333333
$ws.Cells.Item(4, 8).Value = 28.388248
$ws.Cells.Item(4, 9).Value = 0.07254428564686972
$ws.Cells.Item(4, 10).Value = 0.07439525120506714
$ws.Cells.Item(4, 13).Value = 49.446822
$ws.Cells.Item(4, 14).Value = 148.340466
$ws.Cells.Item(4, 15).Value = 0.2944913370783471
$ws.Cells.Item(4, 16).Value = 0.2979770110192673
$ws.Cells.Item(4, 17).Value = 467.902881915952
$ws.Cells.Item(4, 18).Value = 4211.125937243567
$ws.Cells.Item(4, 19).Value = 0.02136366367754021
$ws.Cells.Item(4, 20).Value = 0.02216807458811345
$ws.Cells.Item(5, 7).Value = 9.462749333333333
$ws.Cells.Item(5, 8).Value = 28.388248
$ws.Cells.Item(5, 9).Value = 0.07254428564686972
$ws.Cells.Item(5, 10).Value = 0.07439525120506714
$ws.Cells.Item(5, 13).Value = 5.892385
$ws.Cells.Item(5, 14).Value = 11.78477
$ws.Cells.Item(5, 15).Value = 0.03509338450973445
$ws.Cells.Item(5, 16).Value = 0.02367250578914543
$ws.Cells.Item(5, 17).Value = 55.75816223049333
$ws.Cells.Item(5, 18).Value = 334.54897338296
$ws.Cells.Item(5, 19).Value = 0.002545824510189609
$ws.Cells.Item(5, 20).Value = 0.001761122014836881
$ws.Cells.Item(6, 7).Value = 9.462749333333333
$ws.Cells.Item(6, 8).Value = 28.388248
$ws.Cells.Item(6, 9).Value = 0.07254428564686972
$ws.Cells.Item(6, 10).Value = 0.07439525120506714
$ws.Cells.Item(6, 13).Value = 73.40740966666667
$ws.Cells.Item(6, 14).Value = 220.222229
$ws.Cells.Item(6, 15).Value = 0.4371938448176639
$ws.Cells.Item(6, 16).Value = 0.4423685817288764
$ws.Cells.Item(6, 17).Value = 694.6359168849768
$ws.Cells.Item(6, 18).Value = 6251.723251964791
$ws.Cells.Item(6, 19).Value = 0.03171591516150585
$ws.Cells.Item(6, 20).Value = 0.03291012176294904
$ws.Cells.Item(7, 9).Value = 0.3231336970688258
$ws.Cells.Item(7, 10).Value = 0.3313784449305509
$ws.Cells.Item(7, 13).Value = 3.790335666666667
$ws.Cells.Item(7, 14).Value = 11.371007
$ws.Cells.Item(7, 15).Value = 0.02257417106507704
$ws.Cells.Item(7, 16).Value = 0.02284136466268864
$ws.Cells.Item(7, 17).Value = 159.7622061441722
$ws.Cells.Item(7, 18).Value = 1437.85985529755
$ws.Cells.Item(7, 19).Value = 0.007294475354522459
$ws.Cells.Item(7, 20).Value = 0.007569135902013399
$ws.Cells.Item(8, 9).Value = 0.3231336970688258
$ws.Cells.Item(8, 10).Value = 0.3313784449305509
$ws.Cells.Item(8, 15).Value = 0.2106472625291773
$ws.Cells.Item(8, 16).Value = 0.2131405368000222
$ws.Cells.Item(8, 19).Value = 0.06806722871848062
$ws.Cells.Item(8, 20).Value = 0.07063017963645418
$ws.Cells.Item(9, 9).Value = 0.3231336970688258
$ws.Cells.Item(9, 10).Value = 0.3313784449305509
$ws.Cells.Item(9, 13).Value = 49.446822
$ws.Cells.Item(9, 14).Value = 148.340466
$ws.Cells.Item(9, 15).Value = 0.2944913370783471
$ws.Cells.Item(9, 16).Value = 0.2979770110192673
$ws.Cells.Item(9, 17).Value = 2084.1777785041
$ws.Cells.Item(9, 18).Value = 18757.6000065369
$ws.Cells.Item(9, 19).Value = 0.0951600745048681
$ws.Cells.Item(9, 20).Value = 0.09874315853661843
$ws.Cells.Item(10, 9).Value = 0.3231336970688258
$ws.Cells.Item(10, 10).Value = 0.3313784449305509
$ws.Cells.Item(10, 13).Value = 5.892385
$ws.Cells.Item(10, 14).Value = 11.78477
$ws.Cells.Item(10, 15).Value = 0.03509338450973445
$ws.Cells.Item(10, 16).Value = 0.02367250578914543
$ws.Cells.Item(10, 17).Value = 248.3633403050833
$ws.Cells.Item(10, 18).Value = 1490.1800418305
$ws.Cells.Item(10, 19).Value = 0.01133985507928836
$ws.Cells.Item(10, 20).Value = 0.007844558156016476
$ws.Cells.Item(11, 9).Value = 0.3231336970688258
$ws.Cells.Item(11, 10).Value = 0.3313784449305509
$ws.Cells.Item(11, 13).Value = 73.40740966666667
$ws.Cells.Item(11, 14).Value = 220.222229
$ws.Cells.Item(11, 15).Value = 0.4371938448176639
$ws.Cells.Item(11, 16).Value = 0.4423685817288764
$ws.Cells.Item(11, 17).Value = 3094.113753252206
$ws.Cells.Item(11, 18).Value = 27847.02377926985
$ws.Cells.Item(11, 19).Value = 0.1412720634116663
$ws.Cells.Item(11, 20).Value = 0.1465914126994484
$ws.Cells.Item(12, 7).Value = 32.300192
$ws.Cells.Item(12, 8).Value = 96.900576
$ws.Cells.Item(12, 9).Value = 0.2476229975407503
$ws.Cells.Item(12, 10).Value = 0.2539410918713864
$ws.Cells.Item(12, 13).Value = 3.790335666666667
$ws.Cells.Item(12, 14).Value = 11.371007
$ws.Cells.Item(12, 15).Value = 0.02257417106507704
$ws.Cells.Item(12, 16).Value = 0.02284136466268864
$ws.Cells.Item(12, 17).Value = 122.4285697777813
$ws.Cells.Item(12, 18).Value = 1101.857128000032
$ws.Cells.Item(12, 19).Value = 0.005589883906132049
$ws.Cells.Item(12, 20).Value = 0.005800361082275655
$ws.Cells.Item(13, 7).Value = 32.300192
$ws.Cells.Item(13, 8).Value = 96.900576
$ws.Cells.Item(13, 9).Value = 0.2476229975407503
$ws.Cells.Item(13, 10).Value = 0.2539410918713864
$ws.Cells.Item(13, 15).Value = 0.2106472625291773
$ws.Cells.Item(13, 16).Value = 0.2131405368000222
$ws.Cells.Item(13, 17).Value = 1142.422594597451
$ws.Cells.Item(13, 18).Value = 10281.80335137706
$ws.Cells.Item(13, 19).Value = 0.05216110657122826
$ws.Cells.Item(13, 20).Value = 0.05412514063705104
$ws.Cells.Item(14, 7).Value = 32.300192
$ws.Cells.Item(14, 8).Value = 96.900576
$ws.Cells.Item(14, 9).Value = 0.2476229975407503
$ws.Cells.Item(14, 10).Value = 0.2539410918713864
$ws.Cells.Item(14, 13).Value = 49.446822
$ws.Cells.Item(14, 14).Value = 148.340466
$ws.Cells.Item(14, 15).Value = 0.2944913370783471
$ws.Cells.Item(14, 16).Value = 0.2979770110192673
$ws.Cells.Item(14, 17).Value = 1597.141844389824
$ws.Cells.Item(14, 18).Value = 14374.27659950842
$ws.Cells.Item(14, 19).Value = 0.07292282763712382
$ws.Cells.Item(14, 20).Value = 0.07566860753080488
$ws.Cells.Item(15, 7).Value = 32.300192
$ws.Cells.Item(15, 8).Value = 96.900576
$ws.Cells.Item(15, 9).Value = 0.2476229975407503
$ws.Cells.Item(15, 10).Value = 0.2539410918713864
$ws.Cells.Item(15, 13).Value = 5.892385
$ws.Cells.Item(15, 14).Value = 11.78477
$ws.Cells.Item(15, 15).Value = 0.03509338450973445
$ws.Cells.Item(15, 16).Value = 0.02367250578914543
$ws.Cells.Item(15, 17).Value = 190.32516683792
$ws.Cells.Item(15, 18).Value = 1141.95100102752
$ws.Cells.Item(15, 19).Value = 0.008689929066150577
$ws.Cells.Item(15, 20).Value = 0.006011421967427307
$ws.Cells.Item(16, 7).Value = 32.300192
$ws.Cells.Item(16, 8).Value = 96.900576
$ws.Cells.Item(16, 9).Value = 0.2476229975407503
$ws.Cells.Item(16, 10).Value = 0.2539410918713864
$ws.Cells.Item(16, 13).Value = 73.40740966666667
$ws.Cells.Item(16, 14).Value = 220.222229
$ws.Cells.Item(16, 15).Value = 0.4371938448176639
$ws.Cells.Item(16, 16).Value = 0.4423685817288764
$ws.Cells.Item(16, 17).Value = 2371.07342645599
$ws.Cells.Item(16, 18).Value = 21339.6608381039
$ws.Cells.Item(16, 19).Value = 0.1082592503601156
$ws.Cells.Item(16, 20).Value = 0.1123355606538275
$ws.Cells.Item(17, 7).Value = 9.736177999999999
$ws.Cells.Item(17, 8).Value = 19.472356
$ws.Cells.Item(17, 9).Value = 0.07464047213559308
$ws.Cells.Item(17, 10).Value = 0.0510299478916239
$ws.Cells.Item(17, 13).Value = 3.790335666666667
$ws.Cells.Item(17, 14).Value = 11.371007
$ws.Cells.Item(17, 15).Value = 0.02257417106507704
$ws.Cells.Item(17, 16).Value = 0.02284136466268864
$ws.Cells.Item(17, 17).Value = 36.90338273041533
$ws.Cells.Item(17, 18).Value = 221.420296382492
$ws.Cells.Item(17, 19).Value = 0.001684946786366995
$ws.Cells.Item(17, 20).Value = 0.001165593648510581
$ws.Cells.Item(18, 7).Value = 9.736177999999999
$ws.Cells.Item(18, 8).Value = 19.472356
$ws.Cells.Item(18, 9).Value = 0.07464047213559308
$ws.Cells.Item(18, 10).Value = 0.0510299478916239
$ws.Cells.Item(18, 15).Value = 0.2106472625291773
$ws.Cells.Item(18, 16).Value = 0.2131405368000222
$ws.Cells.Item(18, 17).Value = 344.3580066713726
$ws.Cells.Item(18, 18).Value = 2066.148040028236
$ws.Cells.Item(18, 19).Value = 0.01572281112924802
$ws.Cells.Item(18, 20).Value = 0.01087655048649788
$ws.Cells.Item(19, 7).Value = 9.736177999999999
$ws.Cells.Item(19, 8).Value = 19.472356
$ws.Cells.Item(19, 9).Value = 0.07464047213559308
$ws.Cells.Item(19, 10).Value = 0.0510299478916239
$ws.Cells.Item(19, 13).Value = 49.446822
$ws.Cells.Item(19, 14).Value = 148.340466
$ws.Cells.Item(19, 15).Value = 0.2944913370783471
$ws.Cells.Item(19, 16).Value = 0.2979770110192673
$ws.Cells.Item(19, 17).Value = 481.4230605263159
$ws.Cells.Item(19, 18).Value = 2888.538363157896
$ws.Cells.Item(19, 19).Value = 0.02198097243936992
$ws.Cells.Item(19, 20).Value = 0.01520575134521505
$ws.Cells.Item(20, 7).Value = 9.736177999999999
$ws.Cells.Item(20, 8).Value = 19.472356
$ws.Cells.Item(20, 9).Value = 0.07464047213559308
$ws.Cells.Item(20, 10).Value = 0.0510299478916239
$ws.Cells.Item(20, 13).Value = 5.892385
$ws.Cells.Item(20, 14).Value = 11.78477
$ws.Cells.Item(20, 15).Value = 0.03509338450973445
$ws.Cells.Item(20, 16).Value = 0.02367250578914543
$ws.Cells.Item(20, 17).Value = 57.36930920452999
$ws.Cells.Item(20, 18).Value = 229.47723681812
$ws.Cells.Item(20, 19).Value = 0.002619386788642488
$ws.Cells.Item(20, 20).Value = 0.001208006736884256
$ws.Cells.Item(21, 7).Value = 9.736177999999999
$ws.Cells.Item(21, 8).Value = 19.472356
$ws.Cells.Item(21, 9).Value = 0.07464047213559308
$ws.Cells.Item(21, 10).Value = 0.0510299478916239
$ws.Cells.Item(21, 13).Value = 73.40740966666667
$ws.Cells.Item(21, 14).Value = 220.222229
$ws.Cells.Item(21, 15).Value = 0.4371938448176639
$ws.Cells.Item(21, 16).Value = 0.4423685817288764
$ws.Cells.Item(21, 17).Value = 714.7076070335872
$ws.Cells.Item(21, 18).Value = 4288.245642201524
$ws.Cells.Item(21, 19).Value = 0.03263235499196565
$ws.Cells.Item(21, 20).Value = 0.02257404567451613
$ws.Cells.Item(22, 7).Value = 36.79199966666667
$ws.Cells.Item(22, 8).Value = 110.375999
$ws.Cells.Item(22, 9).Value = 0.2820585476079611
$ws.Cells.Item(22, 10).Value = 0.2892552641013719
$ws.Cells.Item(22, 13).Value = 3.790335666666667
$ws.Cells.Item(22, 14).Value = 11.371007
$ws.Cells.Item(22, 15).Value = 0.02257417106507704
$ws.Cells.Item(22, 16).Value = 0.02284136466268864
$ws.Cells.Item(22, 17).Value = 139.4540285845548
$ws.Cells.Item(22, 18).Value = 1255.086257260993
$ws.Cells.Item(22, 19).Value = 0.006367237904069292
$ws.Cells.Item(22, 20).Value = 0.006606984967941746
$ws.Cells.Item(23, 7).Value = 36.79199966666667
$ws.Cells.Item(23, 8).Value = 110.375999
$ws.Cells.Item(23, 9).Value = 0.2820585476079611
$ws.Cells.Item(23, 10).Value = 0.2892552641013719
$ws.Cells.Item(23, 15).Value = 0.2106472625291773
$ws.Cells.Item(23, 16).Value = 0.2131405368000222
$ws.Cells.Item(23, 17).Value = 1301.292937194363
$ws.Cells.Item(23, 18).Value = 11711.63643474927
$ws.Cells.Item(23, 19).Value = 0.05941486092657265
$ws.Cells.Item(23, 20).Value = 0.06165202226279858
$ws.Cells.Item(24, 7).Value = 36.79199966666667
$ws.Cells.Item(24, 8).Value = 110.375999
$ws.Cells.Item(24, 9).Value = 0.2820585476079611
$ws.Cells.Item(24, 10).Value = 0.2892552641013719
$ws.Cells.Item(24, 13).Value = 49.446822
$ws.Cells.Item(24, 14).Value = 148.340466
$ws.Cells.Item(24, 15).Value = 0.2944913370783471
$ws.Cells.Item(24, 16).Value = 0.2979770110192673
$ws.Cells.Item(24, 17).Value = 1819.247458541726
$ws.Cells.Item(24, 18).Value = 16373.22712687553
$ws.Cells.Item(24, 19).Value = 0.0830637988194451
$ws.Cells.Item(24, 20).Value = 0.08619141901851557
$ws.Cells.Item(25, 7).Value = 36.79199966666667
$ws.Cells.Item(25, 8).Value = 110.375999
$ws.Cells.Item(25, 9).Value = 0.2820585476079611
$ws.Cells.Item(25, 10).Value = 0.2892552641013719
$ws.Cells.Item(25, 13).Value = 5.892385
$ws.Cells.Item(25, 14).Value = 11.78477
$ws.Cells.Item(25, 15).Value = 0.03509338450973445
$ws.Cells.Item(25, 16).Value = 0.02367250578914543
$ws.Cells.Item(25, 17).Value = 216.7926269558717
$ws.Cells.Item(25, 18).Value = 1300.75576173523
$ws.Cells.Item(25, 19).Value = 0.009898389065463418
$ws.Cells.Item(25, 20).Value = 0.006847396913980517
$ws.Cells.Item(26, 7).Value = 36.79199966666667
$ws.Cells.Item(26, 8).Value = 110.375999
$ws.Cells.Item(26, 9).Value = 0.2820585476079611
$ws.Cells.Item(26, 10).Value = 0.2892552641013719
$ws.Cells.Item(26, 13).Value = 73.40740966666667
$ws.Cells.Item(26, 14).Value = 220.222229
$ws.Cells.Item(26, 15).Value = 0.4371938448176639
$ws.Cells.Item(26, 16).Value = 0.4423685817288764
$ws.Cells.Item(26, 17).Value = 2700.805391986863
$ws.Cells.Item(26, 18).Value = 24307.24852788177
$ws.Cells.Item(26, 19).Value = 0.1233142608924106
$ws.Cells.Item(26, 20).Value = 0.1279574409381355
